$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.014.45'
$ws.Range("E2").Value = '  +1.31%  '

$ws.Range("D3").Value = '1.959.23'
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '''245.55'
$ws.Range("E5").Value = '  -1.31%  '

$ws.Range("D6").Value = '''0.9994'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '''0.4877'
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("E8").Value = '  +0.70%  '

$ws.Range("D9").Value = '''0.06831'
$ws.Range("E9").Value = '  +0.61%  '

$ws.Range("E10").Value = '  -1.22%  '

$ws.Range("D11").Value = '''107.07'
$ws.Range("E11").Value = '  -1.96%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '''0.07826'
$ws.Range("E12").Value = '  +0.73%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.940.15'
$ws.Range("E13").Value = '  -1.10%  '

$ws.Range("D14").Value = '''5.498'
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").Value = '''0.7052'
$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("D16").Value = '''284.33'
$ws.Range("E16").Value = '  -3.21%  '

$ws.Range("D17").Value = '31.037.52'
$ws.Range("E17").Value = '  +1.29%  '

$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("D19").Value = '''0.000007711'
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").Value = '2.216.78'
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").Value = '''0.9993'
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").Value = '''5.533'
$ws.Range("E22").Value = '  -1.60%  '

$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").Value = '''6.525'
$ws.Range("E24").Value = '  -1.72%  '

$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("D26").Value = '''168.75'
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").Value = '''20.02'
$ws.Range("E27").Value = '  -0.51%  '

$ws.Range("D28").Value = '''2.222'
$ws.Range("E28").Value = '  +2.17%  '

$ws.Range("D29").Value = '''0.1057'
$ws.Range("E29").Value = '  -0.59%  '

$ws.Range("D30").Value = '''1.388'
$ws.Range("E30").Value = '  -3.39%  '

$ws.Range("D31").Value = '''1.583'
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("D32").Value = '''4.621'
$ws.Range("E32").Value = '  -1.49%  '

$ws.Range("D33").Value = '''4.455'
$ws.Range("E33").Value = '  +0.53%  '

$ws.Range("D34").Value = '''0.04956'
$ws.Range("E34").Value = '  -2.88%  '

$ws.Range("D35").Value = '''0.7639'
$ws.Range("E35").Value = '  -0.80%  '

$ws.Range("D36").Value = '''1.175'
$ws.Range("E36").Value = '  -0.50%  '

$ws.Range("D37").Value = '''2.729'

$ws.Range("D38").Value = '''0.02026'
$ws.Range("E38").Value = '  -0.73%  '

$ws.Range("D39").Value = '''2.703'
$ws.Range("E39").Value = '  -0.58%  '

$ws.Range("D40").Value = '''6.548'
$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").Value = '''78.38'
$ws.Range("E41").Value = '  +11.47%  '

$ws.Range("D42").Value = '''2.128'
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").Value = '''0.9076'
$ws.Range("E43").Value = '  +3.78%  '

$ws.Range("D44").Value = '''0.4488'
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").Value = '''109.08'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").Value = '''8.187'
$ws.Range("E46").Value = '  +8.68%  '

$ws.Range("D47").Value = '''0.9997'
$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("E48").Value = '  +11.70%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.1264'
$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.372'
$ws.Range("E50").Value = '  -0.30%  '

$ws.Range("E51").Value = '  -0.19%  '

Write-Output "done"
